# Weekly fruit/vegetable price update — reshuffle per-row data (columns D, M..T)
# across rows 2-20 (excluding rows 12 and 17, which are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row number -> source row number (values currently present
# at the source row, before the edit, should end up at the target row).
$rowMap = @{
    2  = 18
    3  = 13
    4  = 20
    5  = 10
    6  = 11
    7  = 2
    8  = 14
    9  = 19
    10 = 4
    11 = 3
    13 = 8
    14 = 5
    15 = 6
    16 = 7
    18 = 9
    19 = 15
    20 = 16
}

# Capture the "before" snapshot of every row that participates in the
# reshuffle so that overwriting a row doesn't clobber a value we still
# need to read from it later.
$snapshot = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $snapshot[$srcRow] = @{
            D = $ws.Cells.Item($srcRow, 4).Value2
            M = $ws.Cells.Item($srcRow, 13).Value2
            N = $ws.Cells.Item($srcRow, 14).Value2
            O = $ws.Cells.Item($srcRow, 15).Value2
            P = $ws.Cells.Item($srcRow, 16).Value2
            Q = $ws.Cells.Item($srcRow, 17).Value2
            R = $ws.Cells.Item($srcRow, 18).Value2
            S = $ws.Cells.Item($srcRow, 19).Value2
            T = $ws.Cells.Item($srcRow, 20).Value2
        }
    }
}

foreach ($targetRow in $rowMap.Keys) {
    $srcRow = $rowMap[$targetRow]
    $data = $snapshot[$srcRow]

    $ws.Cells.Item($targetRow, 4).Value  = $data.D   # D - Fecha
    $ws.Cells.Item($targetRow, 13).Value = $data.M   # M - Volumen
    $ws.Cells.Item($targetRow, 14).Value = $data.N   # N - Precio minimo
    $ws.Cells.Item($targetRow, 15).Value = $data.O   # O - Precio maximo
    $ws.Cells.Item($targetRow, 16).Value = $data.P   # P - Precio promedio ponderado
    $ws.Cells.Item($targetRow, 17).Value = $data.Q   # Q - Unidad de comercializacion
    $ws.Cells.Item($targetRow, 18).Value = $data.R   # R - Origen
    $ws.Cells.Item($targetRow, 19).Value = $data.S   # S - Precio $/Kg
    $ws.Cells.Item($targetRow, 20).Value = $data.T   # T - Kg / unidad
}
